# Insert a new "Week_Start_Date" column after "Week" (column A) on the
# "Forecast Comparison" sheet, fix the Week labels to drop the leading
# zero (W01 -> W1 ... W09 -> W9), and populate the new column with the
# week start dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Shift B:I right by one column to make room for the new column B.
$ws.Range("B:B").Insert()

# New header for column B.
$ws.Range("B1").Value = "Week_Start_Date"

# Store the new column as plain text so the dates aren't auto-converted
# into date serial numbers by Excel's cell-value parser.
$ws.Range("B2:B17").NumberFormat = "@"

# Week start dates for rows 2-17 (weeks W1-W16).
$weekStartDates = @(
    "2025-01-05",
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20"
)

for ($i = 0; $i -lt $weekStartDates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $weekStartDates[$i]

    # Fix the Week label in column A: W01 -> W1 ... W09 -> W9 (W10-W16 unchanged).
    $weekLabel = $ws.Cells.Item($row, 1).Value2
    if ($weekLabel -match '^W0(\d)$') {
        $ws.Cells.Item($row, 1).Value = "W" + $matches[1]
    }
}
